$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$ws.Range("B2").Value = 17.52856320514249
$ws.Range("C2").Value = 9.532289039123967
$ws.Range("E2").Value = 11.64646648567682
$ws.Range("F2").Value = 16.86991607391233
$ws.Range("G2").Value = 40.02804150606046
$ws.Range("H2").Value = 16.86384290741739
$ws.Range("L2").Value = 9.914138470075724
$ws.Range("M2").Value = 15.90873117907651
$ws.Range("N2").Value = 19.07096936297475
$ws.Range("B3").Value = 17.12748390953544
$ws.Range("C3").Value = 9.241980144588794
$ws.Range("E3").Value = 11.6677460222269
$ws.Range("F3").Value = 15.89584955866808
$ws.Range("G3").Value = 39.9055691191341
$ws.Range("H3").Value = 16.90347327196797
$ws.Range("L3").Value = 9.923512210824526
$ws.Range("M3").Value = 15.83772752787294
$ws.Range("N3").Value = 19.14528044247534
$ws.Range("B4").Value = 16.88078476152909
$ws.Range("C4").Value = 9.057197283084198
$ws.Range("E4").Value = 11.68156407098756
$ws.Range("F4").Value = 15.26997757108491
$ws.Range("G4").Value = 39.84539046968646
$ws.Range("H4").Value = 16.93166410638931
$ws.Range("L4").Value = 9.930664770324887
$ws.Range("M4").Value = 15.79710713986227
$ws.Range("N4").Value = 19.19294441569158
$ws.Range("B5").Value = 16.78029921566081
$ws.Range("C5").Value = 8.98031723524752
$ws.Range("E5").Value = 11.68738468930925
$ws.Range("F5").Value = 15.00819731993403
$ws.Range("G5").Value = 39.82465272483914
$ws.Range("H5").Value = 16.94411890466057
$ws.Range("L5").Value = 9.933931056623345
$ws.Range("M5").Value = 15.78131473468802
$ws.Range("N5").Value = 19.21288181953291
$ws.Range("B6").Value = 16.76362110384686
$ws.Range("C6").Value = 8.967457920152576
$ws.Range("E6").Value = 11.68836266694763
$ws.Range("F6").Value = 14.96433081551593
$ws.Range("G6").Value = 39.82143801671054
$ws.Range("H6").Value = 16.94624530911301
$ws.Range("L6").Value = 9.9344946607984
$ws.Range("M6").Value = 15.77873870650891
$ws.Range("N6").Value = 19.21622350460605
$ws.Range("B7").Value = 16.87942918163707
$ws.Range("C7").Value = 9.056166759683098
$ws.Range("E7").Value = 11.68164180135354
$ws.Range("F7").Value = 15.26647399323137
$ws.Range("G7").Value = 39.84509545861592
$ws.Range("H7").Value = 16.93182816605605
$ws.Range("L7").Value = 9.930707396851272
$ws.Range("M7").Value = 15.79689106219876
$ws.Range("N7").Value = 19.1932112151687
$ws.Range("B8").Value = 17.39045325709681
$ws.Range("C8").Value = 9.433582211850345
$ws.Range("E8").Value = 11.65364780894364
$ws.Range("F8").Value = 16.53996406344768
$ws.Range("G8").Value = 39.98270257952639
$ws.Range("H8").Value = 16.87670491673887
$ws.Range("L8").Value = 9.917080752715446
$ws.Range("M8").Value = 15.88363911961085
$ws.Range("N8").Value = 19.09617020754662
$ws.Range("B9").Value = 18.38254283209877
$ws.Range("C9").Value = 10.11939669580422
$ws.Range("E9").Value = 11.60469956160679
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 40.37107644699181
$ws.Range("H9").Value = 16.79935457307404
$ws.Range("L9").Value = 9.901430908597995
$ws.Range("M9").Value = 16.07678127359488
$ws.Range("N9").Value = 18.92194917469503
$ws.Range("B10").Value = 19.09696656333255
$ws.Range("C10").Value = 10.58726036944736
$ws.Range("E10").Value = 11.57233328944124
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 40.72741894644192
$ws.Range("H10").Value = 16.76144695345738
$ws.Range("L10").Value = 9.896663118208215
$ws.Range("M10").Value = 16.23188159358402
$ws.Range("N10").Value = 18.80362990212461
$ws.Range("B11").Value = 19.41721940738847
$ws.Range("C11").Value = 10.79177705102379
$ws.Range("E11").Value = 11.55838353077946
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 40.90456284807368
$ws.Range("H11").Value = 16.74834433261519
$ws.Range("L11").Value = 9.895950038455316
$ws.Range("M11").Value = 16.30511103701307
$ws.Range("N11").Value = 18.75187991607354
$ws.Range("B12").Value = 19.53768567675676
$ws.Range("C12").Value = 10.86798999304984
$ws.Range("E12").Value = 11.55321191938282
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 40.97376520162718
$ws.Range("H12").Value = 16.7439806679699
$ws.Range("L12").Value = 9.895888751251674
$ws.Range("M12").Value = 16.33320796472092
$ws.Range("N12").Value = 18.73257987249461
$ws.Range("B13").Value = 19.5117788709627
$ws.Range("C13").Value = 10.85163155906267
$ws.Range("E13").Value = 11.55432079448367
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 40.95876757749433
$ws.Range("H13").Value = 16.74489382826436
$ws.Range("L13").Value = 9.895892677152402
$ws.Range("M13").Value = 16.32714077900868
$ws.Range("N13").Value = 18.73672331975983
$ws.Range("B14").Value = 19.42714708629424
$ws.Range("C14").Value = 10.79807204565984
$ws.Range("E14").Value = 11.55795583976796
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 40.91021389689547
$ws.Range("H14").Value = 16.74797333315189
$ws.Range("L14").Value = 9.895940816567952
$ws.Range("M14").Value = 16.30741534723032
$ws.Range("N14").Value = 18.75028615616893
$ws.Range("B15").Value = 19.37519909276405
$ws.Range("C15").Value = 10.76510368579858
$ws.Range("E15").Value = 11.5601968349065
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 40.88074837918113
$ws.Range("H15").Value = 16.7499375658724
$ws.Range("L15").Value = 9.895997468363481
$ws.Range("M15").Value = 16.29538014101315
$ws.Range("N15").Value = 18.75863235498448
$ws.Range("B16").Value = 19.07593100858859
$ws.Range("C16").Value = 10.57372409151678
$ws.Range("E16").Value = 11.57326047463934
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 40.71614144960797
$ws.Range("H16").Value = 16.76238679687388
$ws.Range("L16").Value = 9.896738964426707
$ws.Range("M16").Value = 16.22714815745159
$ws.Range("N16").Value = 18.80705347563229
$ws.Range("B17").Value = 18.89103369688265
$ws.Range("C17").Value = 10.45416092488574
$ws.Range("E17").Value = 11.58147249168232
$ws.Range("F17").Value = 20.20408069617459
$ws.Range("G17").Value = 40.6189853086923
$ws.Range("H17").Value = 16.77108660794909
$ws.Range("L17").Value = 9.897566310882787
$ws.Range("M17").Value = 16.1859622473787
$ws.Range("N17").Value = 18.83728827028715
$ws.Range("B18").Value = 18.78424620276096
$ws.Range("C18").Value = 10.38461045822634
$ws.Range("E18").Value = 11.58626868688873
$ws.Range("F18").Value = 19.95656407809808
$ws.Range("G18").Value = 40.56452227794237
$ws.Range("H18").Value = 16.77648021966106
$ws.Range("L18").Value = 9.898179247421055
$ws.Range("M18").Value = 16.16252596110195
$ws.Range("N18").Value = 18.85487384398437
$ws.Range("B19").Value = 18.74801825803607
$ws.Range("C19").Value = 10.36092895014085
$ws.Range("E19").Value = 11.58790512270922
$ws.Range("F19").Value = 19.87204792380562
$ws.Range("G19").Value = 40.54632685241696
$ws.Range("H19").Value = 16.7783732587719
$ws.Range("L19").Value = 9.898410336044444
$ws.Range("M19").Value = 16.15463478489242
$ws.Range("N19").Value = 18.86086161506883
$ws.Range("B20").Value = 18.91076274377523
$ws.Range("C20").Value = 10.46696973069844
$ws.Range("E20").Value = 11.58059077071614
$ws.Range("F20").Value = 20.24955283636157
$ws.Range("G20").Value = 40.62918120773604
$ws.Range("H20").Value = 16.77012014635791
$ws.Range("L20").Value = 9.897464056499764
$ws.Range("M20").Value = 16.19032052607326
$ws.Range("N20").Value = 18.8340495223271
$ws.Range("B21").Value = 19.45202833123786
$ws.Range("C21").Value = 10.81383751289847
$ws.Range("E21").Value = 11.55688513344285
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 40.92441805766788
$ws.Range("H21").Value = 16.74705255982033
$ws.Range("L21").Value = 9.895921016871702
$ws.Range("M21").Value = 16.31319938309153
$ws.Range("N21").Value = 18.74629438741951
$ws.Range("B22").Value = 19.80102406071857
$ws.Range("C22").Value = 11.03333408199227
$ws.Range("E22").Value = 11.54203811380797
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 41.12971735703786
$ws.Range("H22").Value = 16.73546319182126
$ws.Range("L22").Value = 9.896128918788204
$ws.Range("M22").Value = 16.3956354960609
$ws.Range("N22").Value = 18.69066909328076
$ws.Range("B23").Value = 19.61523221428834
$ws.Range("C23").Value = 10.91685472866422
$ws.Range("E23").Value = 11.54990327579279
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 41.01903057008941
$ws.Range("H23").Value = 16.74132888841011
$ws.Range("L23").Value = 9.895906878512797
$ws.Range("M23").Value = 16.35144918173031
$ws.Range("N23").Value = 18.72019982617724
$ws.Range("B24").Value = 18.90184475075406
$ws.Range("C24").Value = 10.46118139601125
$ws.Range("E24").Value = 11.58098916313552
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 40.62456729865181
$ws.Range("H24").Value = 16.77055586285295
$ws.Range("L24").Value = 9.897509858043902
$ws.Range("M24").Value = 16.18834939298592
$ws.Range("N24").Value = 18.83551312722608
$ws.Range("B25").Value = 18.11612047183936
$ws.Range("C25").Value = 9.939999144893296
$ws.Range("E25").Value = 11.61730771496726
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 40.25342898237319
$ws.Range("H25").Value = 16.81696908964721
$ws.Range("L25").Value = 9.904481027324135
$ws.Range("M25").Value = 16.02214919648833
$ws.Range("N25").Value = 18.96737191274167
